# Update the "取得日時" (acquisition timestamp) column (A) for rows 2-9
# on the "ランサーズ" sheet, replacing the old timestamp with the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-10 06:38:35"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
